# Fix typos in "01 - What is Vue.js.pptx"

$p = $ppt.ActivePresentation

# --- Slide 1: "What is Vue?" bullet list -------------------------------
# Reword the trailing parenthetical on the last bullet.
$s1 = $p.Slides.Item(1)
$shRect6_s1 = $s1.Shapes.Item(3)                 # "Rectangle 6"
$tr1 = $shRect6_s1.TextFrame.TextRange
$lastPara1 = $tr1.Paragraphs(4, 1)               # 4th paragraph (single paragraph)
$lastRun1 = $lastPara1.Runs($lastPara1.Runs().Count)
$lastRun1.Text = "(Vue-Router, Vuex, etc …)"

# --- Slide 2: "Vue Pros" bullet list -------------------------------------
# Move the bold span boundary: "(" stays with normal text, ")" moves
# onto the bold run by itself.
$s2 = $p.Slides.Item(2)
$shRect6_s2 = $s2.Shapes.Item(3)                 # "Rectangle 6"
$tr2 = $shRect6_s2.TextFrame.TextRange
$firstPara2 = $tr2.Paragraphs(1, 1)              # 1st paragraph ("Small Size ...")
$run1 = $firstPara2.Runs(1)
$run2 = $firstPara2.Runs(2)
$run1.Text = "Small Size (18 KB after gzipping"
$run2.Text = ")"

# Nudge the "What is Vue?" label box up slightly.
$shRect7_s2 = $s2.Shapes.Item(4)                 # "Rectangle 7"
$shRect7_s2.Top = 1247593 / 12700

# --- Slide 3: screenshot picture -----------------------------------------
# Reposition the picture slightly (almost unchanged horizontally, moved
# up vertically).
$s3 = $p.Slides.Item(3)
$shPic_s3 = $s3.Shapes.Item(1)                   # "Picture 8"
$shPic_s3.Left = 2756178 / 12700
$shPic_s3.Top = 1405167 / 12700
